$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-14 09:17:53"
$wsZhCn.Range("H3").Value = "2016-03-14 09:18:36"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-14 09:18:00"
$wsDeDe.Range("H3").Value = "2016-03-14 09:18:49"
